$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.423.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.88%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.435.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.21%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "406.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.21%  "

$ws.Range("E7").Value = "  +2.54%  "

$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("E9").Value = "  +7.56%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.140"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +19.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.56"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.89%  "

$ws.Range("E12").Value = "  +0.64%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.438.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.548.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.52%  "

$ws.Range("E18").Value = "  +2.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000165"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +33.63%  "

$ws.Range("E20").Value = "  +1.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "84.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "315.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.84%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.69%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.94%  "

$ws.Range("E29").Value = "  +9.48%  "

$ws.Range("E30").Value = "  +2.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "44.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.49%  "

$ws.Range("E32").Value = "  +3.44%  "

$ws.Range("E33").Value = "  +3.38%  "

$ws.Range("E34").Value = "  -0.04%  "

$ws.Range("E35").Value = "  +3.95%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.85%  "

$ws.Range("E37").Value = "  +0.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.97"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.320"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +16.71%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "143.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.71%  "

$ws.Range("E42").Value = "  +4.77%  "

$ws.Range("E43").Value = "  +3.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.00%  "

$ws.Range("E45").Value = "  +3.68%  "

$ws.Range("E46").Value = "  +0.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.109.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.32%  "

$ws.Range("E49").Value = "  +12.14%  "

$ws.Range("E50").Value = "  +1.68%  "

$ws.Range("E51").Value = "  +32.68%  "
